# B6-PowerPoint.pptx edit: Mon, Jun 15, 2020 4:04:45 AM
#
# 1) Three tables (the ones that use the deck's custom table style
#    {955AA49C-679A-4681-8C65-C9CD46FED8D2}) get switched to the built-in
#    table style {C7178C74-372F-40D9-971D-4DB362729D77}.
# 2) The design theme's 12-colour scheme is swapped from the "Integral"
#    (Red Violet) palette to the stock "Office Theme" palette.

$p = $ppt.ActivePresentation

# --- 1) Re-style every table that is still on the old custom style ----
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)
    for ($j = 1; $j -le $s.Shapes.Count; $j++) {
        $sh = $s.Shapes.Item($j)
        if ($sh.HasTable) {
            $tbl = $sh.Table
            if ($tbl.Style -eq "{955AA49C-679A-4681-8C65-C9CD46FED8D2}") {
                $tbl.ApplyStyle("{C7178C74-372F-40D9-971D-4DB362729D77}")
            }
        }
    }
}

# --- 2) Swap the theme's colour scheme to the Office Theme palette ----
$slide = $p.Slides.Item(1)
$colors = $slide.ThemeColorScheme

$colors.Colors(1).RGB  = 0         # dk1      000000
$colors.Colors(2).RGB  = 16777215  # lt1      FFFFFF
$colors.Colors(3).RGB  = 6968388   # dk2      44546A
$colors.Colors(4).RGB  = 15132391  # lt2      E7E6E6
$colors.Colors(5).RGB  = 13998939  # accent1  5B9BD5
$colors.Colors(6).RGB  = 3243501   # accent2  ED7D31
$colors.Colors(7).RGB  = 10855845  # accent3  A5A5A5
$colors.Colors(8).RGB  = 49407     # accent4  FFC000
$colors.Colors(9).RGB  = 12874308  # accent5  4472C4
$colors.Colors(10).RGB = 4697456   # accent6  70AD47
$colors.Colors(11).RGB = 12673797  # hlink    0563C1
$colors.Colors(12).RGB = 7491477   # folHlink 954F72
